# Updated cryptos list (prices / 1h volume %) and re-ordered
# Bittensor / InternetComputer(DFINITY) rows to match the refreshed
# coinranking.com data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.994.18'
$ws.Range('E2').Value = '  -2.03%  '
$ws.Range('D3').Value = '2.535.08'
$ws.Range('E3').Value = '  -3.49%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = "'579.24"
$ws.Range('E5').Value = '  -2.88%  '
$ws.Range('D6').Value = "'166.65"
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'0.523"
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').Value = '2.535.08'
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = "'0.350"
$ws.Range('E12').Value = '  -3.41%  '
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('D14').Value = "'26.41"
$ws.Range('E14').Value = '  -4.53%  '
$ws.Range('D15').Value = '2.999.64'
$ws.Range('E15').Value = '  -3.62%  '
$ws.Range('E16').Value = '  -3.23%  '
$ws.Range('D17').Value = '65.836.85'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').Value = '2.544.60'
$ws.Range('E18').Value = '  -3.00%  '
$ws.Range('D19').Value = "'11.26"
$ws.Range('E19').Value = '  -7.01%  '
$ws.Range('D20').Value = "'7.62"
$ws.Range('E20').Value = '  -5.95%  '
$ws.Range('D21').Value = "'345.83"
$ws.Range('E21').Value = '  -2.91%  '
$ws.Range('E22').Value = '  -3.47%  '
$ws.Range('D23').Value = "'4.53"
$ws.Range('E23').Value = '  -3.08%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = "'1.92"
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').Value = "'68.68"
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').Value = "'9.93"
$ws.Range('E27').Value = '  -6.11%  '
$ws.Range('E28').Value = '  -2.81%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').Value = '0.0₃0970'
$ws.Range('E30').Value = '  -3.34%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = "'523.26"
$ws.Range('E31').Value = '  -4.61%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'8.15"
$ws.Range('E32').Value = '  +2.93%  '
$ws.Range('E33').Value = '  -3.41%  '
$ws.Range('E34').Value = '  -4.95%  '
$ws.Range('E35').Value = '  -4.22%  '
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').Value = "'156.65"
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('E38').Value = '  -3.84%  '
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('D40').Value = "'18.26"
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('D43').Value = "'5.04"
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').Value = "'147.24"
$ws.Range('E46').Value = '  -2.84%  '
$ws.Range('D47').Value = '0.0₆0281'
$ws.Range('E47').Value = '  -4.38%  '
$ws.Range('E48').Value = '  -4.06%  '
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('D50').Value = "'1.69"
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('E51').Value = '  -2.10%  '
